# Insert two new weekly price records into the Acelga / Macroferia Regional
# de Talca dataset.
#
# The data block (rows 2-269) holds one row per reporting date; this change
# adds two additional dated observations:
#   - a new row at position 237 (date 44748)
#   - a new row at position 259 (date 44747, counted AFTER the first insert)
# Every row below each insertion point shifts down by one, which is why the
# sheet's used range grows from A1:R269 to A1:R271.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$peso = [char]36
$region = "Regi" + [char]0x00F3 + "n del Maule"
$unidad = $peso + "/docena de atados (4 kilos)"

function Fill-AcelgaRow($row, $fecha, $volumen, $precioMin, $precioMax, $precioProm, $precioKg) {
    $ws.Range("A$row").Value = 5
    $ws.Range("B$row").Value = "Macroferia Regional de Talca"
    $ws.Range("C$row").Value = "Maule"
    $ws.Range("D$row").Value = $fecha
    $ws.Range("E$row").Value = 7
    $ws.Range("F$row").Value = 100112009
    $ws.Range("G$row").Value = "Acelga"
    $ws.Range("H$row").Value = "Sin especificar"
    $ws.Range("I$row").Value = "Primera"
    $ws.Range("J$row").Value = $volumen
    $ws.Range("K$row").Value = $precioMin
    $ws.Range("L$row").Value = $precioMax
    $ws.Range("M$row").Value = $precioProm
    $ws.Range("N$row").Value = $unidad
    $ws.Range("O$row").Value = $region
    $ws.Range("P$row").Value = $precioKg
    $ws.Range("Q$row").Value = 4
    $ws.Range("R$row").Value = "Hortaliza"
}

# First insertion: push old rows 237-269 down to 238-270, then populate the
# freshly-opened row 237.
$ws.Rows.Item(237).Insert()
Fill-AcelgaRow 237 44748 500 2500 2500 2500 625

# Second insertion: row 259 now holds what used to be row 258's data: push
# it (and everything below) down one more, then populate the new row 259.
$ws.Rows.Item(259).Insert()
Fill-AcelgaRow 259 44747 500 3000 3000 3000 750
